$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking text values must be forced to Text format ("@") before
# assignment, otherwise Excel auto-coerces them to numbers (dropping
# formatting such as trailing zeros, e.g. "85.00" -> 85).
$textCells = @{
    'D2' = '70.942.43'
    'E2' = '  -0.23%  '
    'D3' = '3.847.55'
    'E3' = '  +1.12%  '
    'D4' = '0.999'
    'E4' = '  -0.08%  '
    'D5' = '704.16'
    'E5' = '  -0.93%  '
    'D6' = '172.56'
    'E6' = '  -0.19%  '
    'D7' = '3.845.85'
    'E7' = '  +1.10%  '
    'E9' = '  -1.07%  '
    'D11' = '7.34'
    'E11' = '  -1.98%  '
    'E12' = '  -0.75%  '
    'E13' = '  -2.43%  '
    'D14' = '36.42'
    'E14' = '  +0.30%  '
    'D15' = '4.491.77'
    'E15' = '  +1.05%  '
    'D16' = '3.837.34'
    'E16' = '  +0.78%  '
    'D17' = '70.935.38'
    'E17' = '  -0.28%  '
    'D18' = '7.18'
    'E18' = '  -0.83%  '
    'E19' = '  +0.71%  '
    'D20' = '17.37'
    'E20' = '  -3.10%  '
    'D21' = '10.70'
    'E21' = '  -4.49%  '
    'D22' = '493.29'
    'E22' = '  +1.82%  '
    'D23' = '0.715'
    'E23' = '  -0.09%  '
    'D24' = '85.00'
    'E24' = '  +1.33%  '
    'E25' = '  -0.75%  '
    'D26' = '12.14'
    'E26' = '  -2.27%  '
    'D27' = '10.52'
    'E27' = '  -0.49%  '
    'D28' = '2.12'
    'E28' = '  -3.07%  '
    'D29' = '3.18'
    'E29' = '  +0.38%  '
    'D30' = '0.999'
    'E30' = '  -0.06%  '
    'D31' = '7.50'
    'E31' = '  -0.96%  '
    'D32' = '2.27'
    'E32' = '  -1.24%  '
    'E33' = '  +1.75%  '
    'D34' = '29.42'
    'E34' = '  -0.79%  '
    'D35' = '3.799.53'
    'E35' = '  +1.20%  '
    'D36' = '9.15'
    'E36' = '  -1.17%  '
    'D37' = '0.999'
    'E37' = '  -0.30%  '
    'E38' = '  -0.44%  '
    'D39' = '2.38'
    'E39' = '  +6.48%  '
    'D40' = '6.02'
    'E40' = '  +0.46%  '
    'D41' = '1.03'
    'E41' = '  +6.18%  '
    'D42' = '3.31'
    'E42' = '  -6.01%  '
    'E43' = '  +0.05%  '
    'D45' = '164.05'
    'E45' = '  +1.04%  '
    'D46' = '0.000311'
    'E46' = '  -6.31%  '
    'D47' = '48.65'
    'E47' = '  -1.68%  '
    'D48' = '0.299'
    'E48' = '  -0.90%  '
    'D49' = '8.63'
    'E49' = '  +0.52%  '
    'D50' = '43.26'
    'E50' = '  -3.81%  '
    'D51' = '412.30'
    'E51' = '  +3.11%  '
}
foreach ($addr in $textCells.Keys) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $textCells[$addr]
}

# Plain text fields (coin name / link) - no numeric coercion risk.
$ws.Range('B40').Value = 'Filecoin'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('B41').Value = 'Mantle'
$ws.Range('C41').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('B45').Value = 'Monero'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('B46').Value = 'FLOKI'
$ws.Range('C46').Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
